$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force specific Price cells to remain text (they contain plain decimal
# numbers which Excel would otherwise auto-convert to numeric values).
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D12",
    "D13",
    "D14",
    "D19",
    "D20",
    "D21",
    "D22",
    "D24",
    "D26",
    "D27",
    "D28",
    "D30",
    "D31",
    "D32",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D43",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range("D2").Value = "96.158.34"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "3.573.53"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "239.92"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").Value = "653.60"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "1.63"
$ws.Range("E7").Value = "  +11.54%  "
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "1.07"
$ws.Range("E9").Value = "  +7.36%  "
$ws.Range("B10").Value = "USDC"
$ws.Range("C10").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D10").Value = "1.00"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").Value = "3.572.81"
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("D12").Value = "43.41"
$ws.Range("E12").Value = "  +2.30%  "
$ws.Range("D13").Value = "0.202"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").Value = "6.38"
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("D15").Value = "4.238.98"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").Value = "96.008.71"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").Value = "3.576.10"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("D19").Value = "7.77"
$ws.Range("E19").Value = "  -5.11%  "
$ws.Range("D20").Value = "12.64"
$ws.Range("E20").Value = "  -3.46%  "
$ws.Range("D21").Value = "17.77"
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("D22").Value = "0.521"
$ws.Range("E22").Value = "  +8.17%  "
$ws.Range("E23").Value = "  -4.96%  "
$ws.Range("D24").Value = "501.95"
$ws.Range("E24").Value = "  -0.78%  "
$ws.Range("E25").Value = "  +4.97%  "
$ws.Range("D26").Value = "0.0000198"
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("D27").Value = "95.93"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").Value = "12.92"
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("D29").Value = "3.766.90"
$ws.Range("E29").Value = "  -1.29%  "
$ws.Range("D30").Value = "0.153"
$ws.Range("E30").Value = "  +10.43%  "
$ws.Range("D31").Value = "3.00"
$ws.Range("E31").Value = "  -4.51%  "
$ws.Range("D32").Value = "11.35"
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  +2.85%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").Value = "31.33"
$ws.Range("E36").Value = "  -2.17%  "
$ws.Range("D37").Value = "619.29"
$ws.Range("E37").Value = "  +7.84%  "
$ws.Range("D38").Value = "8.80"
$ws.Range("E38").Value = "  +8.27%  "
$ws.Range("D39").Value = "0.563"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("E40").Value = "  +10.78%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").Value = "0.902"
$ws.Range("E43").Value = "  -2.45%  "
$ws.Range("E44").Value = "  +5.68%  "
$ws.Range("D45").Value = "5.71"
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0423"
$ws.Range("E46").Value = "  +2.51%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").Value = "23.52"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").Value = "2.27"
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("D49").Value = "33.36"
$ws.Range("E49").Value = "  -4.85%  "
$ws.Range("D50").Value = "3.49"
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("D51").Value = "8.23"
$ws.Range("E51").Value = "  +2.58%  "
